$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("quality_comparison")

$c1 = $ws1.Range("C1")
$c1.ClearFormats()
$c1.Borders(9).LineStyle = 1
$c1.Borders(8).LineStyle = 1

$d1 = $ws1.Range("D1")
$d1.ClearFormats()
$d1.Borders(9).LineStyle = 1
$d1.Borders(8).LineStyle = 1
$d1.Borders(10).LineStyle = 1

Write-Host "done sheet1"
